# Suzanny.xlsx - "Listas sem duplicação de professores"
# Replace cells that contained Python-list-looking strings (duplicated teacher
# entries) with a simple "-" placeholder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("E12", "E14", "E15", "B18", "E18", "B19", "C19", "E19", "B20", "C20", "D20", "E20", "B21", "C21", "E21", "F21")

foreach ($cell in $cells) {
    $ws.Range($cell).Value = "-"
}
